$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

# The cell that previously showed "Sunny Days" (A11) is re-labelled "Holidays",
# and the header cell C1 (also "Sunny Days") is updated to "Holidays" too.
# This removes the now-unused "Sunny Days" shared string from the workbook.
$ws.Range("A11").Value = "Holidays"
$ws.Range("C1").Value = "Holidays"

# Move / restore the active selection to A12.
$ws.Range("A12").Select()
